$d = $word.ActiveDocument

$p12 = $d.Paragraphs(12).Range
$p12.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Planned</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>storypoints</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>:</w:t></w:r></w:p>')

$p13 = $d.Paragraphs(13).Range
$p13.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>Reviews</w:t></w:r><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>eva</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p14 = $d.Paragraphs(14).Range
$p14.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">Admin can accept reject review </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>eva</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p15 = $d.Paragraphs(15).Range
$p15.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">Other things </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>todo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-AU"/></w:rPr><w:t>:</w:t></w:r></w:p>')

$p16 = $d.Paragraphs(16).Range
$p16.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>file</w:t></w:r></w:p>')

$p17 = $d.Paragraphs(17).Range
$p17.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>next</w:t></w:r><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>year</w:t></w:r><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> login</w:t></w:r></w:p>')

$p18 = $d.Paragraphs(18).Range
$p18.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="0070C0"/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve">Statistics (how many tutors are available, graphics, </w:t></w:r><w:r><w:rPr><w:color w:val="0070C0"/><w:lang w:val="en-AU"/></w:rPr><w:t>which subject, which classes, how many requests…)</w:t></w:r></w:p>')

$p19 = $d.Paragraphs(19).Range
$p19.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>Comment code</w:t></w:r></w:p>')

$p20 = $d.Paragraphs(20).Range
$p20.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t>Todos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="00B050"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> in code</w:t></w:r></w:p>')

$p21 = $d.Paragraphs(21).Range
$p21.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="0070C0"/><w:lang w:val="en-AU"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/><w:lang w:val="en-AU"/></w:rPr><w:t>Verschiedene</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/><w:lang w:val="en-AU"/></w:rPr><w:t>Bilder</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/><w:lang w:val="en-AU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/><w:lang w:val="en-AU"/></w:rPr><w:t>Startseite</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p22 = $d.Paragraphs(22).Range
$p22.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>Cmd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>strg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> taste bei </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>fächer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>auswahl</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p23 = $d.Paragraphs(23).Range
$p23.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="00B050"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve">Oben </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>namen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve"> anzeigen</w:t></w:r></w:p>')

$p24 = $d.Paragraphs(24).Range
$p24.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">Wenn man bestätigt wird, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>email</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> schicken</w:t></w:r></w:p>')

$p25 = $d.Paragraphs(25).Range
$p25.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">Nur </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>jpg</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>bild</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> vorher </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>info</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">, bei </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>fach</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> rote </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>fehlermeldung</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p26 = $d.Paragraphs(26).Range
$p26.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve">Kontaktiere uns </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>name</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t xml:space="preserve"> auslesen</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:rPr><w:color w:val="0070C0"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="0070C0"/></w:rPr><w:t>Mail an alle Tutoren</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr></w:p>')

Write-Output "done"
Write-Output $d.Paragraphs.Count